$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New evaluation results for "Query 2" (row 4) -- x1 and x2 columns.
$ws.Range("D4").Value = 0.74275198707938395
$ws.Range("F4").Value = 0.88397478133732699

# Move / leave the selection where the author left it after entering the new data.
$ws.Range("H4:I4").Select() | Out-Null
